# AHD_DDX41_v4_VAR.docx content update (Release 1.0 template refresh)
#
# Applies the substantive text edits described by the commit:
#   1. Re-word the assay detection-limit exception list so ASXL1
#      c.1934dup;p.Gly646Trpfs*12 is listed first (~5%-10%) and CEBPA is
#      joined by the newly-added TERT gene (~10%).
#   2. Update the VAF / measurement-of-uncertainty (CV%) table values.
#   3. Refresh the "Date authorised" field result to the new issue date.

$d = $word.ActiveDocument

# 1) Detection-limit exception list: reorder genes and fold TERT into the
#    10% group, widening the ASXL1 indel limit to a 5%-10% range.
$d.Content.Find.Execute(
    "CEBPA (detection limit ~ 10%) and ASXL1 c.1934dup;p.Gly646Trpfs*12 (detection limit ~ 5%)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ASXL1 c.1934dup;p.Gly646Trpfs*12 (detection limit ~ 5%-10%), CEBPA and TERT (detection limit ~ 10%)",
    2
) | Out-Null

# 2) Measurement-of-uncertainty CV% table: new VAF bands and averages.
$d.Content.Find.Execute(
    "VAFs of 5%, 10%-20%, 30%-40% and 50% are on average, 10.2%, 10.4%, 3.5% and 4.4%, respectively.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "VAFs of 2%, 5%-10%, 20%-40% and 50% are on average, 15.4%, 8.6%, 4.0% and 1.8%, respectively.",
    2
) | Out-Null

# 3) Date authorised field result.
$d.Content.Find.Execute(
    "30-Oct-2023",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "16-Nov-2023",
    2
) | Out-Null
